# Auto-generated Excel COM-interop edit script
# Updates market price / profit data cells across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# matching the scheduled-runner data refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1322.0521
$ws.Range("I15").Value = 1322.0521
$ws.Range("K15").Value = 3966.156300000001
$ws.Range("M15").Value = -3797.156300000001
# Row 28
$ws.Range("H28").Value = 7774.067
$ws.Range("I28").Value = 217.52632
$ws.Range("J28").Value = 20826.273
$ws.Range("K28").Value = 217.52632
$ws.Range("L28").Value = 20826.273
$ws.Range("M28").Value = 267.47368
$ws.Range("N28").Value = -21796.273
# Row 64
$ws.Range("H64").Value = 3930.353
$ws.Range("I64").Value = 3763.5557
$ws.Range("J64").Value = 3990.4
$ws.Range("K64").Value = 3763.5557
$ws.Range("L64").Value = 3990.4
$ws.Range("M64").Value = -3515.5557
$ws.Range("N64").Value = -4486.4
# Row 67
$ws.Range("H67").Value = 3930.353
$ws.Range("I67").Value = 3763.5557
$ws.Range("J67").Value = 3990.4
$ws.Range("K67").Value = 3763.5557
$ws.Range("L67").Value = 3990.4
$ws.Range("M67").Value = -2905.5557
$ws.Range("N67").Value = -5706.4
# Row 137
$ws.Range("H137").Value = 1136.9333
$ws.Range("I137").Value = 941.3182
$ws.Range("J137").Value = 1674.875
$ws.Range("K137").Value = 2823.9546
$ws.Range("L137").Value = 5024.625
$ws.Range("M137").Value = -273.9546
$ws.Range("N137").Value = -10124.625

$ws = $wb.Worksheets.Item("ARM")
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("N47").Value = 0
$ws.Range("L47").ClearContents()
# Row 63
$ws.Range("H63").Value = 1284382.2
$ws.Range("I63").Value = 2085146.6
$ws.Range("J63").Value = 3159.3
$ws.Range("K63").Value = 2085146.6
$ws.Range("L63").Value = 3159.3
$ws.Range("M63").Value = -2084460.6
$ws.Range("N63").Value = -4531.3
# Row 66
$ws.Range("H66").Value = 1284382.2
$ws.Range("I66").Value = 2085146.6
$ws.Range("J66").Value = 3159.3
$ws.Range("K66").Value = 10425733
$ws.Range("L66").Value = 15796.5
$ws.Range("M66").Value = -10422301
$ws.Range("N66").Value = -22660.5
# Row 74
$ws.Range("H74").Value = 1164.25
$ws.Range("I74").Value = 966.6667
$ws.Range("J74").Value = 1282.8
$ws.Range("K74").Value = 966.6667
$ws.Range("L74").Value = 1282.8
$ws.Range("M74").Value = -92.66669999999999
$ws.Range("N74").Value = -3030.8
# Row 77
$ws.Range("H77").Value = 1164.25
$ws.Range("I77").Value = 966.6667
$ws.Range("J77").Value = 1282.8
$ws.Range("K77").Value = 4833.3335
$ws.Range("L77").Value = 6414
$ws.Range("M77").Value = -465.3334999999997
$ws.Range("N77").Value = -15150
# Row 122
$ws.Range("H122").Value = 38463680
$ws.Range("I122").Value = 47621120
$ws.Range("J122").Value = 2440
$ws.Range("K122").Value = 142863360
$ws.Range("L122").Value = 7320
$ws.Range("M122").Value = -142860910
$ws.Range("N122").Value = -12220

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2333.3157
$ws.Range("I105").Value = 1955.5555
$ws.Range("J105").Value = 2673.3
$ws.Range("K105").Value = 1955.5555
$ws.Range("L105").Value = 2673.3
$ws.Range("M105").Value = -208.5554999999999
$ws.Range("N105").Value = -6167.3
# Row 107
$ws.Range("H107").Value = 819
$ws.Range("I107").Value = 784.7222
$ws.Range("J107").Value = 973.25
$ws.Range("K107").Value = 784.7222
$ws.Range("L107").Value = 973.25
$ws.Range("M107").Value = 1135.2778
$ws.Range("N107").Value = -4813.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6851.403
$ws.Range("I31").Value = 2664.1667
$ws.Range("J31").Value = 27787.584
$ws.Range("K31").Value = 2664.1667
$ws.Range("L31").Value = 27787.584
$ws.Range("M31").Value = -2369.1667
$ws.Range("N31").Value = -28377.584
# Row 34
$ws.Range("H34").Value = 6851.403
$ws.Range("I34").Value = 2664.1667
$ws.Range("J34").Value = 27787.584
$ws.Range("K34").Value = 2664.1667
$ws.Range("L34").Value = 27787.584
$ws.Range("M34").Value = -2462.1667
$ws.Range("N34").Value = -28191.584
# Row 62
$ws.Range("H62").Value = 7563.75
$ws.Range("J62").Value = 4500
$ws.Range("L62").Value = 4500
$ws.Range("N62").Value = -5748
# Row 65
$ws.Range("H65").Value = 7563.75
$ws.Range("J65").Value = 4500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -28740
# Row 69
$ws.Range("H69").Value = 20000
$ws.Range("I69").Value = 20000
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = -19251
$ws.Range("K69").ClearContents()
$ws.Range("M69").ClearContents()
# Row 72
$ws.Range("H72").Value = 20000
$ws.Range("I72").Value = 20000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 60000
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = -56256
$ws.Range("M72").ClearContents()
# Row 141
$ws.Range("H141").Value = 78214.836
$ws.Range("J141").Value = 78214.836
$ws.Range("L141").Value = 78214.836
$ws.Range("N141").Value = -88574.836

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 60.347828
$ws.Range("I23").Value = 58.3
$ws.Range("J23").Value = 61.923077
$ws.Range("K23").Value = 174.9
$ws.Range("L23").Value = 185.769231
$ws.Range("M23").Value = 60.10000000000002
$ws.Range("N23").Value = -655.769231
# Row 80
$ws.Range("H80").Value = 5046.6665
$ws.Range("I80").Value = 3800
$ws.Range("J80").Value = 5135.7144
$ws.Range("K80").Value = 11400
$ws.Range("L80").Value = 15407.1432
$ws.Range("M80").Value = -10464
$ws.Range("N80").Value = -17279.1432
# Row 83
$ws.Range("H83").Value = 5046.6665
$ws.Range("I83").Value = 3800
$ws.Range("J83").Value = 5135.7144
$ws.Range("K83").Value = 34200
$ws.Range("L83").Value = 46221.4296
$ws.Range("M83").Value = -29520
$ws.Range("N83").Value = -55581.4296

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 34636.668
$ws.Range("I70").Value = 53840.5
$ws.Range("J70").Value = 5092.3076
$ws.Range("K70").Value = 53840.5
$ws.Range("L70").Value = 5092.3076
$ws.Range("M70").Value = -53570.5
$ws.Range("N70").Value = -5632.3076
# Row 73
$ws.Range("H73").Value = 34636.668
$ws.Range("I73").Value = 53840.5
$ws.Range("J73").Value = 5092.3076
$ws.Range("K73").Value = 53840.5
$ws.Range("L73").Value = 5092.3076
$ws.Range("M73").Value = -52904.5
$ws.Range("N73").Value = -6964.3076
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()
# Row 107
$ws.Range("H107").Value = 128.85715
$ws.Range("I107").Value = 110.4
$ws.Range("J107").Value = 175
$ws.Range("K107").Value = 110.4
$ws.Range("L107").Value = 175
$ws.Range("M107").Value = 1809.6
$ws.Range("N107").Value = -4015
# Row 109
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32080
# Row 122
$ws.Range("H122").Value = 37039836
$ws.Range("I122").Value = 62501428
$ws.Range("K122").Value = 187504284
$ws.Range("M122").Value = -187501834

$ws = $wb.Worksheets.Item("WVR")
# Row 25
$ws.Range("H25").Value = 6581.6665
$ws.Range("I25").Value = 6019
$ws.Range("J25").Value = 6694.2
$ws.Range("K25").Value = 6019
$ws.Range("L25").Value = 6694.2
$ws.Range("M25").Value = -5726
$ws.Range("N25").Value = -7280.2
# Row 58
$ws.Range("H58").Value = 5000000
$ws.Range("I58").Value = 5000000
$ws.Range("K58").Value = 5000000
$ws.Range("M58").Value = -4999692
